$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 2500
$ws.Range("I64").Value = 2500
$ws.Range("K64").Value = 2500
$ws.Range("M64").Value = -2252
# Row 67
$ws.Range("H67").Value = 2500
$ws.Range("I67").Value = 2500
$ws.Range("K67").Value = 2500
$ws.Range("M67").Value = -1642
# Row 111
$ws.Range("H111").Value = 5400
$ws.Range("J111").Value = 6750
$ws.Range("L111").Value = 20250
$ws.Range("N111").Value = -26384
# Row 141
$ws.Range("H141").Value = 33333.332
$ws.Range("I141").Value = 15000
$ws.Range("K141").Value = 45000
$ws.Range("M141").Value = -39820

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# Row 63
$ws.Range("H63").Value = 1002.5
$ws.Range("I63").Value = 1002.5
$ws.Range("K63").Value = 1002.5
$ws.Range("M63").Value = -316.5
# Row 66
$ws.Range("H66").Value = 1002.5
$ws.Range("I66").Value = 1002.5
$ws.Range("K66").Value = 5012.5
$ws.Range("M66").Value = -1580.5
# Row 88
$ws.Range("H88").Value = 4718.3335
$ws.Range("I88").Value = 2005
$ws.Range("J88").Value = 6075
$ws.Range("K88").Value = 2005
$ws.Range("L88").Value = 6075
$ws.Range("M88").Value = -1599
$ws.Range("N88").Value = -6887
# Row 91
$ws.Range("H91").Value = 4718.3335
$ws.Range("I91").Value = 2005
$ws.Range("J91").Value = 6075
$ws.Range("K91").Value = 2005
$ws.Range("L91").Value = 6075
$ws.Range("M91").Value = -601
$ws.Range("N91").Value = -8883
# Row 97
$ws.Range("H97").Value = 900
$ws.Range("I97").Value = 900
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -404
$ws.Range("N97").ClearContents()
# Row 112
$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 20000
$ws.Range("N112").Value = -22954
# Row 133
$ws.Range("H133").Value = 69992
$ws.Range("I133").Value = 40000
$ws.Range("J133").Value = 99984
$ws.Range("K133").Value = 40000
$ws.Range("L133").Value = 99984
$ws.Range("M133").Value = -37470
$ws.Range("N133").Value = -105044

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2085.1333
$ws.Range("I20").Value = 1970.7778
$ws.Range("J20").Value = 2256.6667
$ws.Range("K20").Value = 1970.7778
$ws.Range("L20").Value = 2256.6667
$ws.Range("M20").Value = -1723.7778
$ws.Range("N20").Value = -2750.6667
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
# Row 86
$ws.Range("H86").Value = 35169
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 41802.8
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 41802.8
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -44048.8
# Row 89
$ws.Range("H89").Value = 35169
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 41802.8
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 209014
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -220246
# Row 99
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 137
$ws.Range("H137").Value = 99995
$ws.Range("J137").Value = 99995
$ws.Range("L137").Value = 99995
$ws.Range("N137").Value = -110195

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 2500
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 22500
$ws.Range("N132").Value = -27560

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3730
$ws.Range("I80").Value = 3730
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3730
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2732
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 3730
$ws.Range("I83").Value = 3730
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 18650
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -13658
$ws.Range("N83").ClearContents()
# Row 97
$ws.Range("H97").Value = 1911.875
$ws.Range("I97").Value = 1911.875
$ws.Range("K97").Value = 1911.875
$ws.Range("M97").Value = -1415.875
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 10000
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10224
# Row 7
$ws.Range("H7").Value = 4467.9287
$ws.Range("I7").Value = 4427
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 4427
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -4315
$ws.Range("N7").Value = -5224
# Row 15
$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10340
# Row 20
$ws.Range("H20").Value = 13285.714
$ws.Range("J20").Value = 13285.714
$ws.Range("L20").Value = 13285.714
$ws.Range("N20").Value = -13737.714
# Row 24
$ws.Range("H24").Value = 18000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 18000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 18000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -18686
# Row 40
$ws.Range("H40").Value = 7234.4614
$ws.Range("J40").Value = 3713.8572
$ws.Range("L40").Value = 3713.8572
$ws.Range("N40").Value = -3985.8572
# Row 61
$ws.Range("H61").Value = 4250
$ws.Range("I61").Value = 4250
$ws.Range("K61").Value = 4250
$ws.Range("M61").Value = -4048
# Row 82
$ws.Range("H82").Value = 1928
$ws.Range("J82").Value = 2299.6
$ws.Range("L82").Value = 2299.6
$ws.Range("N82").Value = -3021.6
# Row 85
$ws.Range("H85").Value = 1928
$ws.Range("J85").Value = 2299.6
$ws.Range("L85").Value = 2299.6
$ws.Range("N85").Value = -4795.6
# Row 113
$ws.Range("H113").Value = 4250
$ws.Range("I113").Value = 4250
$ws.Range("K113").Value = 4250
$ws.Range("M113").Value = -2080
# Row 122
$ws.Range("H122").Value = 6742.7144
$ws.Range("I122").Value = 7489.8
$ws.Range("K122").Value = 22469.4
$ws.Range("M122").Value = -20019.4
# Row 126
$ws.Range("H126").Value = 4467.9287
$ws.Range("I126").Value = 4427
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 13281
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -10811
$ws.Range("N126").Value = -19940
# Row 136
$ws.Range("H136").Value = 4303.0557
$ws.Range("I136").Value = 4026.7646
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 12080.2938
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -9530.293799999999
$ws.Range("N136").Value = -32100

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4667
$ws.Range("I81").Value = 5200.4
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 10400.8
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -9339.799999999999
$ws.Range("N81").Value = -6122
# Row 84
$ws.Range("H84").Value = 4667
$ws.Range("I84").Value = 5200.4
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 52004
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -46700
$ws.Range("N84").Value = -30608
# Row 126
$ws.Range("H126").Value = 7668
$ws.Range("I126").Value = 7668
$ws.Range("K126").Value = 23004
$ws.Range("M126").Value = -20534
